# Generate Report for Handoff
# The "b1e47540-b6d9-473a-9f28-054d6d17d239.md" file moves from "In Translation"
# to "Ready for handoff" status, and gets a fresh Latest Handoff Datetime stamp
# on the two language sheets (zh-cn, de-de). The Overview sheet's rollup status
# for that file mirrors the same change.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: roll-up status for b1e47540-...md (row 3) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: status + handoff datetime for b1e47540-...md (row 3) ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("D3").Value = "2016-03-10 16:15:07"

# --- de-de sheet: status + handoff datetime for b1e47540-...md (row 3) ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Ready for handoff"
$de.Range("D3").Value = "2016-03-10 16:15:11"
